# Update files from solar experiments
# Adds a new "Solar" worksheet (after "Blad1") containing three PWM / single-charge
# distance experiments plus their averages, and re-points the active sheet/selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Add the new "Solar" sheet right after "Blad1" and make it the active tab.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Solar"

# ---------------------------------------------------------------------------
# 2. Populate "Solar" with the three experiment blocks.
# ---------------------------------------------------------------------------

# --- Block 1: PWM 70 ---
$ws2.Range("A2").Value = "PWM: 70"
$ws2.Range("B2").Value = "Run time:"

$ws2.Range("A3").Value = "Distance: 75cm max"
$ws2.Range("B3").Value = 3
$ws2.Range("D3").Value = "dist single charge"

$ws2.Range("A4").Value = 69.5
$ws2.Range("B4").Formula = "=A4/B3"
$ws2.Range("D4").Value = "21cm"

$ws2.Range("A5").Value = 69.5
$ws2.Range("B5").Formula = "=A5/B3"
$ws2.Range("D5").Value = 28.5
$ws2.Range("E5").Value = 55.7

$ws2.Range("A6").Value = 70.9
$ws2.Range("B6").Formula = "=A6/B3"
$ws2.Range("D6").Value = 27.6
$ws2.Range("E6").Value = 54.3

$ws2.Range("A7").Value = 70
$ws2.Range("B7").Formula = "=A7/B3"

$ws2.Range("A8").Value = 69.9
$ws2.Range("B8").Formula = "=A8/B3"

$ws2.Range("A9").Value = "AVG:"
$ws2.Range("B9").Formula = "=AVERAGE(B4:B8)"

# --- Block 2: PWM 50 ---
$ws2.Range("A11").Value = "PWM: 50"
$ws2.Range("B11").Value = "Run time:"

$ws2.Range("A12").Value = "Distance: 75cm max"
$ws2.Range("B12").Value = 3.5
$ws2.Range("D12").Value = "distance single charge"

$ws2.Range("A13").Value = 73.6
$ws2.Range("B13").Formula = "=A13/B12"
$ws2.Range("D13").Value = 25.4
$ws2.Range("E13").Value = 52.2

$ws2.Range("A14").Value = 71.8
$ws2.Range("B14").Formula = "=A14/B12"
$ws2.Range("D14").Value = 24.8
$ws2.Range("E14").Value = 48.9

$ws2.Range("A15").Value = 72.4
$ws2.Range("B15").Formula = "=A15/B12"
$ws2.Range("D15").Value = 25.4
$ws2.Range("E15").Value = 50.2

$ws2.Range("A16").Value = 72.8
$ws2.Range("B16").Formula = "=A16/B12"
$ws2.Range("D16").Value = 24.3
$ws2.Range("E16").Value = 50.6

$ws2.Range("A17").Value = 72.6
$ws2.Range("B17").Formula = "=A17/B12"
$ws2.Range("D17").Value = 26.1
$ws2.Range("E17").Value = 52.8

$ws2.Range("A18").Value = "AVG:"
$ws2.Range("B18").Formula = "=AVERAGE(B13:B17)"

# --- Block 3: PWM 30 ---
$ws2.Range("A20").Value = "PWM: 30"
$ws2.Range("B20").Value = "Run time:"

$ws2.Range("A21").Value = "Distance: 75cm max"
$ws2.Range("B21").Value = 4.5
$ws2.Range("D21").Value = "distance single charge"

$ws2.Range("A22").Value = 73.2
$ws2.Range("B22").Formula = "=A22/B21"
$ws2.Range("D22").Value = 22
$ws2.Range("E22").Value = 45.2

$ws2.Range("A23").Value = 69.4
$ws2.Range("B23").Formula = "=A23/B21"
$ws2.Range("D23").Value = 20.7
$ws2.Range("E23").Value = 41.4

$ws2.Range("A24").Value = 73.4
$ws2.Range("B24").Formula = "=A24/B21"
$ws2.Range("D24").Value = 21.2
$ws2.Range("E24").Value = 44.5

$ws2.Range("A25").Value = 70.5
$ws2.Range("B25").Formula = "=A25/B21"

$ws2.Range("A26").Value = 74.1
$ws2.Range("B26").Formula = "=A26/B21"
$ws2.Range("D26").Value = 21.7
$ws2.Range("E26").Value = 45.5

$ws2.Range("A27").Value = "AVG:"
$ws2.Range("B27").Formula = "=AVERAGE(B22:B26)"

# ---------------------------------------------------------------------------
# 3. Row heights / column widths on "Solar" (match the source spreadsheet).
# ---------------------------------------------------------------------------
foreach ($r in 1..9) { $ws2.Rows.Item($r).RowHeight = 13.8 }
foreach ($r in 11..18) { $ws2.Rows.Item($r).RowHeight = 13.8 }
foreach ($r in 20..27) { $ws2.Rows.Item($r).RowHeight = 13.8 }

$ws2.Columns.Item(1).ColumnWidth = 18.202380952380967
$ws2.Columns.Item(2).ColumnWidth = 10.687074829931966

# ---------------------------------------------------------------------------
# 4. Column widths on "Blad1" shrink slightly (re-saved alongside the new sheet).
# ---------------------------------------------------------------------------
$ws1.Columns.Item(1).ColumnWidth = 19.68707482993197
$ws1.Columns.Item(2).ColumnWidth = 15.232993197278867
$ws1.Columns.Item(3).ColumnWidth = 15.906462585033966
$ws1.Columns.Item(4).ColumnWidth = 15.906462585033966
$ws1.Columns.Item(5).ColumnWidth = 15.498299319727865
$ws1.Columns.Item(6).ColumnWidth = 15.906462585033966
$ws1.Columns.Item(7).ColumnWidth = 16.447278911564666
$ws1.Columns.Item(8).ColumnWidth = 16.447278911564666
$ws1.Columns.Item(9).ColumnWidth = 16.447278911564666
$ws1.Columns.Item(10).ColumnWidth = 21.17176870748297
$ws1.Columns.Item(11).ColumnWidth = 8.212585034013607
$ws1.Columns.Item(12).ColumnWidth = 19.819727891156468

# ---------------------------------------------------------------------------
# 5. Selection / active sheet bookkeeping.
#    ("Blad1" keeps a pending selection at B34 but "Solar" ends up the active,
#    front-most tab with D12 selected - matches the source workbook.)
# ---------------------------------------------------------------------------
$ws1.Range("B34").Select() | Out-Null
$ws2.Range("D12").Select() | Out-Null
